# RESTdesign.xlsx update:
# - added PUT and PATCH endpoint for Player (/players/{id} row gains PUT + DELETE columns,
#   the /players row which incorrectly had them gets reset to "x")
# - endpoint text for "/games?result{WHITE/BLACK/DRAW}" corrected to include "="
# - status ("Gotowe"/done) column updated for the affected rows
# - minor view/column-width cosmetic tweaks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the query-string endpoint text (missing "=" before the placeholder)
$ws.Range("D5").Value = "/games?result={WHITE/BLACK/DRAW}"

# /games/{gid} (row 6) is now fully done
$ws.Range("J6").Value = "done"

# /players (row 8, LP=6) no longer individually implements PUT/DELETE/PATCH -> reset to "x"
$ws.Range("G8").Value = "x"
$ws.Range("H8").Value = "x"
$ws.Range("I8").Value = "x"
$ws.Range("J8").Value = "done"
# row 8 no longer needs the taller wrapped height now that the long text is gone
$ws.Rows.Item(8).AutoFit()

# /players/{id} (row 9, LP=7) gains the new PUT and DELETE(PATCH) endpoints
$ws.Range("G9").Value = "update players data"
$ws.Range("H9").Value = "delete player"
# match the plain centered/wrapped look used elsewhere in the table (no border)
$ws.Range("E8").Copy()
$ws.Range("G9:H9").PasteSpecial(-4122)
$ws.Range("J9").Value = "GET/PUT/PATCH"

# Column D is a touch wider to fit the updated text
$ws.Columns.Item(4).ColumnWidth = 33.498697916666664

# Selection cosmetic change left over from editing
$ws.Range("K10").Select()
